# Shopee sample workbook update
# - Drops the old "% Giam gia" (discount) computed values in column D (rows 4-7)
#   along with their red-font style (fontId=2 / cellXfs index 2), since that
#   style is no longer referenced by any cell after this edit.
# - Adds a handful of structurally-empty placeholder cells (inlineStr with no
#   text) that the data-export script leaves behind for columns B..G on rows
#   that don't have a value for every field.
# - Appends four new product rows (8-11) with fresh sample data.
# - The used range grows from A1:J7 to A1:G11.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Strip the discount column's old values + the red-font style (D4:D7).
#    ClearFormats() drops the xfId=2 formatting (so that style/font becomes
#    unused and is dropped from the style table on save); Value = "" drops
#    the old negative-percentage numbers.
# ---------------------------------------------------------------------------
$ws.Range("D4:D7").ClearFormats()
$ws.Range("D4:D7").Value = ""

# ---------------------------------------------------------------------------
# 2) Existing rows 2,3,5: add the new blank placeholder cells introduced by
#    the refreshed export (empty inline strings in previously-untouched
#    columns).
# ---------------------------------------------------------------------------
$ws.Range("B2:G2").Value = ""
$ws.Range("D3:E3").Value = ""
$ws.Range("E4").Value = ""
$ws.Range("B5:C5").Value = ""
$ws.Range("E5").Value = ""
$ws.Range("E6").Value = ""
$ws.Range("E7").Value = ""

# ---------------------------------------------------------------------------
# 3) New product rows 8-11.
# ---------------------------------------------------------------------------
$ws.Cells.Item(8, 1).Value = "https://shopee.vn/Samsung-Galaxy-S24-Ultra-12GB-256GB-i.88201679.23626487486"
$ws.Cells.Item(8, 2).Value = "Đen"
$ws.Cells.Item(8, 3).Value = "256GB"
$ws.Cells.Item(8, 4).Value = ""
$ws.Cells.Item(8, 5).Value = ""
$ws.Cells.Item(8, 6).Value = 28990000
$ws.Cells.Item(8, 7).Value = 29500000

$ws.Cells.Item(9, 1).Value = "https://shopee.vn/Google-Pixel-8-Pro-128GB-Ch%C3%ADnh-H%C3%A3ng-i.88201679.23548769421"
$ws.Cells.Item(9, 2).Value = "Xanh"
$ws.Cells.Item(9, 3).Value = "128GB"
$ws.Cells.Item(9, 4).Value = ""
$ws.Cells.Item(9, 5).Value = ""
$ws.Cells.Item(9, 6).Value = 21490000
$ws.Cells.Item(9, 7).Value = 21990000

$ws.Cells.Item(10, 1).Value = "https://shopee.vn/Laptop-Dell-XPS-13-Plus-9320-i7-1260P-16GB-512GB-Windows-11-i.88201679.21845126497"
$ws.Cells.Item(10, 2).Value = ""
$ws.Cells.Item(10, 3).Value = "512GB"
$ws.Cells.Item(10, 4).Value = ""
$ws.Cells.Item(10, 5).Value = ""
$ws.Cells.Item(10, 6).Value = 38990000
$ws.Cells.Item(10, 7).Value = 39900000

$ws.Cells.Item(11, 1).Value = "https://shopee.vn/Laptop-Gaming-Asus-ROG-Strix-G16-G614JV-N4086W-i7-13650HX-16GB-512GB-RTX-4060-Windows-11-i.88201679.22274851269"
$ws.Cells.Item(11, 2).Value = ""
$ws.Cells.Item(11, 3).Value = ""
$ws.Cells.Item(11, 4).Value = ""
$ws.Cells.Item(11, 5).Value = ""
$ws.Cells.Item(11, 6).Value = 39490000
$ws.Cells.Item(11, 7).Value = 39990000
